$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.095691442489624
$ws.Range("B1").Value = 2.186473369598389
$ws.Range("C1").Value = 9.671916961669922
$ws.Range("D1").Value = 1.187099695205688
$ws.Range("E1").Value = 1.236798048019409
